$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3660, 3833, 4203, 4203, 4203, 4203, 4290, 4290, 4290, 4294, 4332, 4374, 4619, 4767)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
